$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": add a new row for the handed-off file
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Duplicate formatting of row 2 into a brand new row 3 (keeps cell styles,
# e.g. the HyperLink style and the date number format, intact).
$wsOverview.Rows("2:2").Copy()
$wsOverview.Rows("3:3").Insert(-4121)

$wsOverview.Range("A3").Value = "69a8c385-e78f-41e3-9726-21a591c205a9ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$wsOverview.Range("B3").Value = "e2e\69a8c385-e78f-41e3-9726-21a591c205a9ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-24 22:28:38"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0c9b79d4c7f87b3d48ecff6e1e52f451e933bcd/e2e/69a8c385-e78f-41e3-9726-21a591c205a9ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md",
    [Type]::Missing,
    [Type]::Missing,
    "e2e\69a8c385-e78f-41e3-9726-21a591c205a9ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
) | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

$wsOverview.Columns("E:F").ColumnWidth = 16.3

# ---------------------------------------------------------------------------
# Sheet "zh-cn": add a new row for the handed-off file
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Rows("2:2").Copy()
$wsZhCn.Rows("3:3").Insert(-4121)

$wsZhCn.Range("A3").Value = "69a8c385-e78f-41e3-9726-21a591c205a9ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "69a8c385-e78f-41e3-9726-21a591c205a9oooooooooooooooooooooooooooooooooooooooo.98db37504a5839b31e178b4fd33b8bb5cff1bbb0.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-24 22:28:33"
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0c9b79d4c7f87b3d48ecff6e1e52f451e933bcd/e2e/69a8c385-e78f-41e3-9726-21a591c205a9ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md",
    [Type]::Missing,
    [Type]::Missing,
    "69a8c385-e78f-41e3-9726-21a591c205a9ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
) | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

$wsZhCn.Columns("C:C").ColumnWidth = 16.3

# ---------------------------------------------------------------------------
# Sheet "de-de": add a new row for the handed-off file
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Rows("2:2").Copy()
$wsDeDe.Rows("3:3").Insert(-4121)

$wsDeDe.Range("A3").Value = "69a8c385-e78f-41e3-9726-21a591c205a9ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "69a8c385-e78f-41e3-9726-21a591c205a9oooooooooooooooooooooooooooooooooooooooo.98db37504a5839b31e178b4fd33b8bb5cff1bbb0.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-24 22:28:38"
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0c9b79d4c7f87b3d48ecff6e1e52f451e933bcd/e2e/69a8c385-e78f-41e3-9726-21a591c205a9ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md",
    [Type]::Missing,
    [Type]::Missing,
    "69a8c385-e78f-41e3-9726-21a591c205a9ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
) | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))

$wsDeDe.Columns("C:C").ColumnWidth = 16.3

Write-Host "Report generated for handoff."
